$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("w01_100.5_optimization_results")

$data = @(
    @(1, 628, 2, 1, -1, 2.965601718751368, 0.001),
    @(2, 637, 0.8948364207292779, 0.2265038621548684, -1, 3, 0.09759263431911268),
    @(3, 627, 1.245351364593948, 1, -0.2582568596808386, 1.272021389224577, 0.001),
    @(4, 638, 2, 0.001, -1, 2.67683929251567, 0.3904905328228748),
    @(5, 636, 1.642464116218913, 0.9202976414121006, -0.1868952056471511, 1.97192326356339, 0.5062142622216756),
    @(6, 628, 1.632362306672351, 1, -1, 1.86579608944121, 1),
    @(7, 637, 1.571086636207806, 0.529333268663417, -1, 2.014440803429157, 0.001),
    @(8, 628, 0.8441961990733206, 1, -1, 1.848961237817112, 1),
    @(9, 627, 0.4033350215703653, 0.6328511528626054, -0.3283095666231202, 1.121984034631293, 0.4355925982014164),
    @(10, 643, 1.049253329906997, 0.1245359487357942, -0.1583826196617993, 1.174234501073174, 0.1984033032562611)
)

$startRow = 12
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

$wb.Save()
